$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ExtremeFlows")
$ws2 = $wb.Worksheets.Item("BasinWaterAccounts")

# Update Extreme Low Flow Method descriptions in column A
$ws1.Range("A8").Value = "       10-year"
$ws1.Range("A9").Value = "       4-year"
$ws1.Range("A10").Value = "       3-year"
$ws1.Range("A5").Value = "D. 85%, 65%, and 50% of 2000 to 2018 average flow"
$ws1.Range("A7").Value = "F. Lowest consecutive flows in Reclamation's ensembles and traces"

# Fix typo "Maximim (maf)" -> "Maximum (maf)" in header E1
$ws1.Range("E1").Value = "Maximum (maf)"

# Widen column E slightly to fit new header text
$ws1.Columns.Item(5).ColumnWidth = 9.81640625

# Update the active selection on the ExtremeFlows sheet
$ws1.Activate()
$ws1.Range("E2").Select()
